# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for
# rows 2-51 of Sheet1, matching the values published by the scraping job.
#
# Several "Price" values look like plain numbers (e.g. 576.75) but must stay
# stored as text, exactly like the rest of the sheet (t="inlineStr"/shared
# string, no numeric coercion, no cell-format change). Assigning such a
# string straight to .Value would make Excel auto-convert it to a Number,
# and pre-seeding NumberFormat="@" would leave a new (visible) style behind.
# The robust fix is the classic Excel trick: prefix with an apostrophe to
# force text entry, then ClearFormats() to drop the quotePrefix style Excel
# attaches for that apostrophe, so the cell ends up with the same (default)
# style it started with.
function Set-TextCell($sheet, $addr, $value) {
    $sheet.Range($addr).Value = "'" + $value
    $sheet.Range($addr).ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.860.59'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '2.445.11'
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("E4").Value = '  -0.19%  '
Set-TextCell $ws "D5" '576.75'
$ws.Range("E5").Value = '  +1.47%  '
Set-TextCell $ws "D6" '145.78'
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '2.442.96'
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("E10").Value = '  +2.44%  '
Set-TextCell $ws "D11" '0.164'
$ws.Range("E11").Value = '  +2.71%  '
Set-TextCell $ws "D12" '5.27'
$ws.Range("E12").Value = '  +1.27%  '
Set-TextCell $ws "D13" '0.354'
$ws.Range("E13").Value = '  +1.92%  '
Set-TextCell $ws "D14" '28.40'
$ws.Range("E14").Value = '  +6.96%  '
Set-TextCell $ws "D15" '0.0000180'
$ws.Range("E15").Value = '  +4.86%  '
$ws.Range("D16").Value = '2.888.34'
$ws.Range("E16").Value = '  +1.67%  '
$ws.Range("D17").Value = '62.787.97'
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("D18").Value = '2.444.53'
$ws.Range("E18").Value = '  +1.13%  '
Set-TextCell $ws "D19" '7.94'
$ws.Range("E19").Value = '  -1.92%  '
Set-TextCell $ws "D20" '11.03'
$ws.Range("E20").Value = '  +2.84%  '
Set-TextCell $ws "D21" '330.97'
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("E23").Value = '  +6.04%  '
$ws.Range("E24").Value = '  +0.03%  '
Set-TextCell $ws "D25" '66.38'
$ws.Range("E25").Value = '  +1.67%  '
Set-TextCell $ws "D26" '649.52'
$ws.Range("E26").Value = '  +9.74%  '
$ws.Range("E27").Value = '  +17.80%  '
$ws.Range("E28").Value = '  +2.38%  '
$ws.Range("D29").Value = '0.0₃0992'
$ws.Range("E29").Value = '  +4.29%  '
$ws.Range("E30").Value = '  +1.96%  '
$ws.Range("E31").Value = '  +1.55%  '
$ws.Range("E32").Value = '  +5.80%  '
$ws.Range("D33").Value = '0.0₆0468'
$ws.Range("E33").Value = '  +62.48%  '
$ws.Range("E34").Value = '  +2.83%  '
$ws.Range("E35").Value = '  +3.13%  '
$ws.Range("E36").Value = '  +0.70%  '
Set-TextCell $ws "D37" '0.998'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("E38").Value = '  +2.73%  '
Set-TextCell $ws "D39" '5.53'
$ws.Range("E39").Value = '  +5.22%  '
$ws.Range("E40").Value = '  +0.07%  '
Set-TextCell $ws "D41" '152.38'
$ws.Range("E41").Value = '  -1.13%  '
Set-TextCell $ws "D42" '18.78'
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("E43").Value = '  +8.18%  '
Set-TextCell $ws "D44" '1.75'
$ws.Range("E44").Value = '  +3.45%  '
Set-TextCell $ws "D45" '42.46'
$ws.Range("E45").Value = '  +2.04%  '
$ws.Range("E46").Value = '  +0.00%  '
Set-TextCell $ws "D47" '14.98'
$ws.Range("E47").Value = '  +27.48%  '
Set-TextCell $ws "D48" '145.34'
$ws.Range("E48").Value = '  +2.08%  '
Set-TextCell $ws "D49" '3.64'
$ws.Range("E49").Value = '  +2.76%  '
Set-TextCell $ws "D50" '20.59'
$ws.Range("E50").Value = '  +4.63%  '
$ws.Range("E51").Value = '  +2.04%  '
